$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (borders, bold, alignment) from the existing header
# cell H1 onto the two new header cells so they match the rest of the header row.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for the new columns I (I0) and J (IF), one row per data row (rows 2-41)
$data = @(
    @(6, 7),
    @(4, 4),
    @(2, 2),
    @(6, 6),
    @(2, 3),
    @(10, 10),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(6, 6),
    @(9, 9),
    @(6, 6),
    @(5, 5),
    @(5, 5),
    @(4, 4),
    @(8, 8),
    @(5, 5),
    @(8, 8),
    @(6, 6),
    @(6, 6),
    @(7, 7),
    @(7, 7),
    @(7, 7),
    @(7, 7),
    @(7, 7),
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(7, 7),
    @(7, 7),
    @(7, 7),
    @(5, 5),
    @(8, 8),
    @(7, 7),
    @(5, 5),
    @(5, 5),
    @(5, 5)
)

for ($idx = 0; $idx -lt $data.Count; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $data[$idx][0]
    $ws.Cells.Item($row, 10).Value = $data[$idx][1]
}
